# Asset production year, insurance page table
# - Remove the "Status" column (old column C)
# - Remove the trailing "Current Cost" / "Useful Life Month" / "Accumulate Depre" /
#   "Net Book Value" columns (old columns Q:T)
# - Insert a new "Production Year" column right after "Sn Engine"
# - Clear the stray styled cell at C2 (row 2 only keeps A2/B2 after the column shuffle)
# - Update the sheet view zoom/selection to match the new layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Status" column entirely (shifts everything after it left by one).
$ws.Range("C1").EntireColumn.Delete() | Out-Null

# Remove the four trailing columns (now at P:S, after the shift above) that are no
# longer part of the table: "Current Cost", "Useful Life Month", "Accumulate Depre",
# "Net Book Value".
$ws.Range("P1:S1").EntireColumn.Delete() | Out-Null

# Insert a new blank column before the (now) "PO No." column so "Production Year"
# can be placed right after "Sn Engine".
$ws.Range("I1").EntireColumn.Insert() | Out-Null

# Set the header text for the newly inserted column (it inherits the bold/shaded
# header style from its neighbors automatically).
$ws.Range("I1").Value = "Production Year"
$ws.Range("I1").EntireColumn.ColumnWidth = 14.33

# Clear the now out-of-place styled cell left over in row 2.
$ws.Range("C2").Clear() | Out-Null

# Match the updated view state.
$ws.Application.ActiveWindow.Zoom = 96
$ws.Range("G23").Select()
